# Updated AoC Days 1-12 plus runtimes to date
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Day 11 (row 15): Part 1 time updated, Part 2 now derived from Total - Part 1,
# and Total becomes a plain recorded value instead of the shared SUM formula.
$ws.Range("B15").Value = 0.011525499983690599
$ws.Range("E15").Value = 0.18739129998721099
$ws.Range("C15").Formula = "=+E15-B15"

# Move the active selection to D15 to match the saved view state.
$ws.Range("D15").Select()
